$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 10.19202271113354
$ws.Range("E2").Value = 13.99772134997901
$ws.Range("F2").Value = 28.00002925415833
$ws.Range("G2").Value = 27.99262492241991
$ws.Range("H2").Value = 13.21923212800039
$ws.Range("J2").Value = 9.644818864116694
$ws.Range("O2").Value = 20.18627902080746

$ws.Range("D3").Value = 10.18792776010369
$ws.Range("E3").Value = 13.96069026936262
$ws.Range("F3").Value = 27.72584099711404
$ws.Range("G3").Value = 27.34581590149835
$ws.Range("H3").Value = 13.19068893700766
$ws.Range("J3").Value = 9.64650525613669
$ws.Range("O3").Value = 20.02620259791254

$ws.Range("D4").Value = 10.18718855004999
$ws.Range("E4").Value = 13.94095794876884
$ws.Range("F4").Value = 27.56376897813321
$ws.Range("G4").Value = 26.94969477667469
$ws.Range("H4").Value = 13.17586455310845
$ws.Range("J4").Value = 9.649197674905501
$ws.Range("O4").Value = 19.93251056327977

$ws.Range("D5").Value = 10.18733472167911
$ws.Range("E5").Value = 13.93367796296429
$ws.Range("F5").Value = 27.49937375959709
$ws.Range("G5").Value = 26.78879404821963
$ws.Range("H5").Value = 13.17050689977642
$ws.Range("J5").Value = 9.650711919510524
$ws.Range("O5").Value = 19.89552268896704

$ws.Range("D6").Value = 10.18738604175381
$ws.Range("E6").Value = 13.93251524037363
$ws.Range("F6").Value = 27.488782674454
$ws.Range("G6").Value = 26.76211597885139
$ws.Range("H6").Value = 13.16965865087502
$ws.Range("J6").Value = 9.650988556930281
$ws.Range("O6").Value = 19.8894539280562

$ws.Range("D7").Value = 10.18718870862639
$ws.Range("E7").Value = 13.94085668000986
$ws.Range("F7").Value = 27.56289374781946
$ws.Range("G7").Value = 26.94752234594239
$ws.Range("H7").Value = 13.17578952584707
$ws.Range("J7").Value = 9.649216407453324
$ws.Range("O7").Value = 19.93200685707587

$ws.Range("D8").Value = 10.19024303432426
$ws.Range("E8").Value = 13.98433288131061
$ws.Range("F8").Value = 27.90422929772433
$ws.Range("G8").Value = 27.76954980634316
$ws.Range("H8").Value = 13.20883245916028
$ws.Range("J8").Value = 9.645056714851281
$ws.Range("O8").Value = 20.13015363980882

$ws.Range("D9").Value = 10.21026356897379
$ws.Range("E9").Value = 14.09315794686801
$ws.Range("F9").Value = 28.62013929065481
$ws.Range("G9").Value = 29.37840297970902
$ws.Range("H9").Value = 13.29485483016439
$ws.Range("J9").Value = 9.650023774114825
$ws.Range("O9").Value = 20.55345977065712

$ws.Range("D10").Value = 10.23344268846872
$ws.Range("E10").Value = 14.187064937627
$ws.Range("F10").Value = 29.16991955979585
$ws.Range("G10").Value = 30.54415329235389
$ws.Range("H10").Value = 13.37067982789178
$ws.Range("J10").Value = 9.661635036204371
$ws.Range("O10").Value = 20.88320280852614

$ws.Range("D11").Value = 10.24580427003785
$ws.Range("E11").Value = 14.23271206016573
$ws.Range("F11").Value = 29.42418676710616
$ws.Range("G11").Value = 31.06815545768677
$ws.Range("H11").Value = 13.40783664559329
$ws.Range("J11").Value = 9.668634509463791
$ws.Range("O11").Value = 21.03673035878794

$ws.Range("D12").Value = 10.25074430124094
$ws.Range("E12").Value = 14.25040869153862
$ws.Range("F12").Value = 29.52098477406992
$ws.Range("G12").Value = 31.26546135151542
$ws.Range("H12").Value = 13.42228238153745
$ws.Range("J12").Value = 9.6715306881252
$ws.Range("O12").Value = 21.09532649437561

$ws.Range("D13").Value = 10.2496688993483
$ws.Range("E13").Value = 14.2465792901557
$ws.Range("F13").Value = 29.50011618201164
$ws.Range("G13").Value = 31.22302119215389
$ws.Range("H13").Value = 13.4191546736729
$ws.Range("J13").Value = 9.670896042929165
$ws.Range("O13").Value = 21.08268715005315

$ws.Range("D14").Value = 10.24620551355723
$ws.Range("E14").Value = 14.23415979427643
$ws.Range("F14").Value = 29.43214062403226
$ws.Range("G14").Value = 31.08441160933743
$ws.Range("H14").Value = 13.40901763072708
$ws.Range("J14").Value = 9.668867865316694
$ws.Range("O14").Value = 21.04154218700499

$ws.Range("D15").Value = 10.24411774134097
$ws.Range("E15").Value = 14.22660571722196
$ws.Range("F15").Value = 29.39056779534323
$ws.Range("G15").Value = 30.99935692190981
$ws.Range("H15").Value = 13.40285702873039
$ws.Range("J15").Value = 9.667657493565203
$ws.Range("O15").Value = 21.01639798227277

$ws.Range("D16").Value = 10.2326712023642
$ws.Range("E16").Value = 14.18413987375149
$ws.Range("F16").Value = 29.15337906325777
$ws.Range("G16").Value = 30.50976287739128
$ws.Range("H16").Value = 13.36830448962617
$ws.Range("J16").Value = 9.661212053365919
$ws.Range("O16").Value = 20.87323624657433

$ws.Range("D17").Value = 10.22611296865166
$ws.Range("E17").Value = 14.15883154443305
$ws.Range("F17").Value = 29.00887749313399
$ws.Range("G17").Value = 30.20763995635958
$ws.Range("H17").Value = 13.3477844908563
$ws.Range("J17").Value = 9.657696969332493
$ws.Range("O17").Value = 20.7862799768707

$ws.Range("D18").Value = 10.22251205291288
$ws.Range("E18").Value = 14.14455113750935
$ws.Range("F18").Value = 28.92616169188486
$ws.Range("G18").Value = 30.03328550724606
$ws.Range("H18").Value = 13.33623314540065
$ws.Range("J18").Value = 9.655836931840598
$ws.Range("O18").Value = 20.73659984591636

$ws.Range("D19").Value = 10.22132231765966
$ws.Range("E19").Value = 14.13976378266693
$ws.Range("F19").Value = 28.89822649347433
$ws.Range("G19").Value = 29.97415906526011
$ws.Range("H19").Value = 13.33236543278872
$ws.Range("J19").Value = 9.65523497299476
$ws.Range("O19").Value = 20.71983797588703

$ws.Range("D20").Value = 10.22679340089652
$ws.Range("E20").Value = 14.1614971365498
$ws.Range("F20").Value = 29.02421940351396
$ws.Range("G20").Value = 30.23986323467775
$ws.Range("H20").Value = 13.34994293450564
$ws.Range("J20").Value = 9.658054425068773
$ws.Range("O20").Value = 20.79550231663106

$ws.Range("D21").Value = 10.24721578547861
$ws.Range("E21").Value = 14.237796628501
$ws.Range("F21").Value = 29.45209347683171
$ws.Range("G21").Value = 31.12515671042001
$ws.Range("H21").Value = 13.41198500623133
$ws.Range("J21").Value = 9.669456935443996
$ws.Range("O21").Value = 21.05361540318239

$ws.Range("D22").Value = 10.26207120930414
$ws.Range("E22").Value = 14.29005360826748
$ws.Range("F22").Value = 29.73468108382642
$ws.Range("G22").Value = 31.69711573386243
$ws.Range("H22").Value = 13.45471638674812
$ws.Range("J22").Value = 9.678340111134014
$ws.Range("O22").Value = 21.22495736481478

$ws.Range("D23").Value = 10.25400542471186
$ws.Range("E23").Value = 14.26194779235944
$ws.Range("F23").Value = 29.58361782623352
$ws.Range("G23").Value = 31.39252465664361
$ws.Range("H23").Value = 13.4317127771795
$ws.Range("J23").Value = 9.673468545754158
$ws.Range("O23").Value = 21.1332825925894

$ws.Range("D24").Value = 10.22648524933168
$ws.Range("E24").Value = 14.1602911815593
$ws.Range("F24").Value = 29.01728219750573
$ws.Range("G24").Value = 30.22529712948349
$ws.Range("H24").Value = 13.34896633560857
$ws.Range("J24").Value = 9.657892318257414
$ws.Range("O24").Value = 20.79133192062998

$ws.Range("D25").Value = 10.20335309967254
$ws.Range("E25").Value = 14.061234597731
$ws.Range("F25").Value = 28.42194760795604
$ws.Range("G25").Value = 28.94502986164284
$ws.Range("H25").Value = 13.26934136391781
$ws.Range("J25").Value = 9.647278360295559
$ws.Range("O25").Value = 20.43547031540906
